$d = $word.ActiveDocument
$cr = [char]13
$lf = [char]10

# ---------------------------------------------------------------------
# The last paragraph in the main body (style "Section Headnote") holds
# one giant run of concatenated text that actually belongs to several
# different resources/sections. Split it into the proper sequence of
# styled paragraphs (Resource/Section Number, Title, Headnote, Case
# Text, ...), matching the structured "casebook" export layout.
# ---------------------------------------------------------------------

# Placeholder markers stand in for paragraphs whose final text needs an
# embedded literal newline char - those are patched up individually
# afterwards (a single Range.Text assignment that both splits on CR and
# embeds a raw LF does not behave the same as doing each separately).
$parts = @(
    "What is a corporation?",
    "1.1",
    "Case of the District Number 1",
    "__PH_RESOURCE_HEADNOTE_1__",
    "This is the body of case 1.",
    "1.2",
    "Case of the District Number 2",
    "__PH_RESOURCE_HEADNOTE_2__",
    "__PH_CASE_TEXT_2__",
    "2",
    "Section Two",
    "__PH_SECTION_HEADNOTE_2__"
)

$target = $d.Paragraphs.Last
$baseIdx = $d.Paragraphs.Count
$joined = [string]::Join($cr, $parts)
$target.Range.Text = $joined

# Styles to apply, in order, to the freshly split paragraphs.
$styles = @(
    "SectionHeadnote",
    "ResourceNumber",
    "ResourceTitle",
    "ResourceHeadnote",
    "CaseText",
    "ResourceNumber",
    "ResourceTitle",
    "ResourceHeadnote",
    "CaseText",
    "SectionNumber",
    "SectionTitle",
    "SectionHeadnote"
)

for ($i = 0; $i -lt $styles.Length; $i++) {
    $d.Paragraphs.Item($baseIdx + $i).Style = $styles[$i]
}

# Bookmarks that wrap the resource/section number runs.
$bookmarks = @{
    1  = "_auto_toc_2"
    5  = "_auto_toc_3"
    9  = "_auto_toc_4"
}
foreach ($offset in $bookmarks.Keys) {
    $bp = $d.Paragraphs.Item($baseIdx + $offset)
    $br = $bp.Range
    $bookRange = $d.Range($br.Start, $br.End)
    $d.Bookmarks.Add($bookmarks[$offset], $bookRange)
}

# Patch up the paragraphs that need an embedded literal newline char at
# the end of their text (done as an isolated Range.Text assignment so
# the LF is kept as real text instead of being read as a paragraph
# break).
$fixups = @{
    3  = "This is the body of case 1."  # unchanged, placeholder replaced below instead
}

$d.Paragraphs.Item($baseIdx + 3).Range.Text = "" + $lf
$d.Paragraphs.Item($baseIdx + 7).Range.Text = "This is an annotatable resource in the casebook." + $lf
$d.Paragraphs.Item($baseIdx + 8).Range.Text = "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;" + $lf
$d.Paragraphs.Item($baseIdx + 11).Range.Text = "This is the second chapter of the casebook." + $lf
